# Updates to Transp sector - fix passenger rail
#
# 1) PCiCDTdtTDM sheet: the B-column (and one C-column) formulas were
#    doubling the Calcs-sheet percentages (e.g. "=Calcs!B5*2"). Remove the
#    stray "*2" so the sheet reports the true (un-doubled) percent change.
# 2) About sheet: remove the stray note "We allow for twice the potential
#    identified in the BLUE Shifts scenario." (the rows below shift up).
# 3) Reset stale cell selections left over from editing.

$wb = $excel.ActiveWorkbook

# --- 1) Fix the doubled passenger-rail (and related) formulas -------------
$calcSheet = $wb.Worksheets.Item("PCiCDTdtTDM")

$calcSheet.Range("B2").Formula = "=Calcs!B5"
$calcSheet.Range("B3").Formula = "=Calcs!C5"
$calcSheet.Range("C3").Formula = "=Calcs!B11"
$calcSheet.Range("B4").Formula = "=Calcs!D5"
$calcSheet.Range("B5").Formula = "=Calcs!E5"
$calcSheet.Range("C5").Formula = "=Calcs!C11"
$calcSheet.Range("B6").Formula = "=Calcs!F5"
$calcSheet.Range("B7").Formula = "=Calcs!G5"

# --- 2) Remove the obsolete note from the About sheet ----------------------
$aboutSheet = $wb.Worksheets.Item("About")
$aboutSheet.Rows("18:19").Delete()

# --- 3) Reset selections back to the top-left cell on the touched sheets --
$aboutSheet.Activate()
$aboutSheet.Range("A1").Select()

$calcSheet.Activate()
$calcSheet.Range("A1").Select()

$aboutSheet.Activate()
